$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 706, pushing existing rows 706..747 down to 707..748
$ws.Rows.Item(706).Insert()

# Populate the newly inserted row with the new data point:
# 2026/01/27, 火, 8:00, rank 25
$ws.Range("A706").Value = "'2026/01/27"
$ws.Range("B706").Value = [char]0x706B
$ws.Range("C706").Value = 8
$ws.Range("D706").Value = 25

# The text assignment above nudges Excel into giving A706 a "text" number
# format/style; restore it to match the plain (unstyled) data rows below it.
$ws.Range("A706").Style = $ws.Range("A707").Style
